# Avance del registro Informe de Visita Verificacion jefe credito
#
# Repurpose the single-row "contrasena/Observacion/Cod cliente" sheet into
# the Informe de Visita Verificacion header/value table (8 columns x 2 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Propagate the existing highlighted/text-formatted style (currently on
# B2:C2) across the rest of the new value row (D2:H2) before we start
# writing values, so the new cells pick up the same yellow fill + "@" text
# number format (and so the new strings are stored as shared strings, not
# numbers).
$ws.Range("B2").Copy($ws.Range("C2:H2")) | Out-Null

# Header row (row 1) - default styling.
$ws.Range("A1").Value2 = "Cod cliente"
$ws.Range("A2").Value2 = "2240270"
$ws.Range("B1").Value2 = "HoraInicioVisita"
$ws.Range("B2").Value2 = "9"
$ws.Range("C2").Value2 = "00"
$ws.Range("C1").Value2 = "MinutosInicioVisita"
$ws.Range("D1").Value2 = "HoraFinVisita"
$ws.Range("E1").Value2 = "MinutosFinVisita"
$ws.Range("D2").Value2 = "12"
$ws.Range("F2").Value2 = "3"
$ws.Range("F1").Value2 = "HoraEmpVisita"
$ws.Range("G1").Value2 = "MinutosEmpVisita"
$ws.Range("H1").Value2 = "Detalle"
$ws.Range("H2").Value2 = "Conforme"
$ws.Range("E2").Value2 = "00"
$ws.Range("G2").Value2 = "00"

# The old third row (previously A3:C3) is no longer part of the table.
$ws.Rows(3).Delete() | Out-Null

# Column C is now a narrow "minutos" column instead of the old wide column.
$ws.Columns.Item(3).ColumnWidth = 6.29

# Move/restore the active selection like the source workbook.
$ws.Range("F10").Select() | Out-Null
